$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the predicted Winner / 2nd Place / 3rd Place picks for the first
# event (row 2, "100m - M"), keeping the existing cell formatting intact.
$ws.Range("B2:D2").ClearContents()

# Move the active selection to D2 (matches the saved selection state in the file)
$ws.Range("D2").Select()
